# Generate and save output file after processing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the old "arts" column (R), shifting everything
# from R onward three columns to the right (R,S,T,... -> U,V,W,...).
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells (row 2) for the inserted columns - numeric zero, matching
# the other general_college_subjects.* columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Lower-case the "Unknown" placeholder text in row 2 (D2:J2).
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"

# The general_college_subjects.arts column (now shifted from R to U) keeps
# its numeric value.
$ws.Range("U2").Value = 0

# sat_act_required.* (now shifted from S,T,U to V,W,X) switch from numeric 0
# to real booleans (FALSE).
$ws.Range("V2").Value = $false
$ws.Range("W2").Value = $false
$ws.Range("X2").Value = $false

# residency_acceptance.out-of-state / .international (now shifted from W,X
# to Z,AA) keep their numeric percentages.
$ws.Range("Z2").Value = 0.28
$ws.Range("AA2").Value = 0.365
